$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0218 to SCD0013
$ws.Name = "SCD0013"

# Update the TC_ID value in B2 from "DGS-233" to "SCD0013-003"
$ws.Range("B2").Value = "SCD0013-003"

# Widen column B to fit the new, longer TC_ID text (best-fit for "SCD0013-003")
$ws.Columns("B").ColumnWidth = 11.67

# Update the view: scroll back to A1 (no frozen/topLeft offset) and move the
# active selection to B3
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B3").Select()
